# Update the "rijm" template styles: reduce the rhyme text font size
# from 24pt (sz 48) to 18pt (sz 36).
#
# - "rijm" is the paragraph style; it previously had no direct run
#   formatting (rPr) of its own, so a new rPr with an 18pt size is added.
# - "rijmChar" is the linked character style; its existing sz is
#   lowered from 48 (24pt) to 36 (18pt), leaving szCs untouched.

$d = $word.ActiveDocument

$rijmStyle = $d.Styles("rijm")
$rijmCharStyle = $d.Styles("rijmChar")

$rijmStyle.Font.Size = 18
$rijmCharStyle.Font.Size = 18
